$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Total price" row (18) and the thin spacer row below it (19) need to move
# down by three rows to make room for three new BOM line items. Using
# Range.Insert (xlShiftDown) on just columns A:C reproduces exactly how this
# was done in Excel (keeps the special "thick border" row flags intact,
# without marking the rows as having an explicit/custom height).
$ws.Range("A18:C18").Insert(-4121)
$ws.Range("A18:C18").Insert(-4121)
$ws.Range("A18:C18").Insert(-4121)

# Three of the old blank spacer rows that used to pad the bottom of the sheet
# are no longer needed (the sheet grew by only one row overall, not three),
# so remove three rows just below the relocated total/spacer rows.
$ws.Range("A23:C25").Delete(-4162)

# New BOM rows 17-19.
$ws.Range("A17").Value = "Square M3 nuts"
$ws.Range("B17").Value = 30

$ws.Range("A18").Value = "M3 nuts"
$ws.Range("B18").Value = 30

$ws.Range("A19").Value = "Various M3 bolts"
$ws.Range("B19").Value = 60

# Rows 19 and 20 keep the same (empty, currency-styled) look as the other
# price cells above them.
$ws.Range("C16").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)

# Selection moved to C17 in the saved file.
[void]$ws.Range("C17").Select()
